$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Test Data")

# Add new rows to Sheet1 ("book list")
$ws1.Cells.Item(5, 1).Value = "7980000029"
$ws1.Cells.Item(5, 2).Value = "test162537"
$ws1.Cells.Item(5, 3).Value = "test162537@gmail.com"
$ws1.Cells.Item(5, 4).Value = "SoftSuave121907"

$ws1.Cells.Item(6, 1).Value = "7980000030"
$ws1.Cells.Item(6, 2).Value = "test162537"
$ws1.Cells.Item(6, 3).Value = "test162537@gmail.com"
$ws1.Cells.Item(6, 4).Value = "SoftSuave121907"

$ws1.Cells.Item(7, 1).Value = "7980000031"
$ws1.Cells.Item(7, 2).Value = "test162537"
$ws1.Cells.Item(7, 3).Value = "test162537@gmail.com"
$ws1.Cells.Item(7, 4).Value = "SoftSuave121907"

# Mark used status in Test Data sheet for rows 30-32 (mobile numbers 7980000029-31)
$ws2.Cells.Item(30, 2).Value = "used"
$ws2.Cells.Item(31, 2).Value = "used"
$ws2.Cells.Item(32, 2).Value = "used"
